$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    "24+55=79",
    "45+41=86",
    "40+16=56",
    "22-1=21",
    "78-51=27",
    "22+73=95",
    "93-1=92",
    "89-87=2",
    "27-19=8",
    "1+63=64",
    "12-7=5",
    "79-16=63",
    "76-60=16",
    "15+34=49",
    "96-15=81",
    "93-49=44",
    "40+44=84",
    "22+32=54",
    "1+55=56",
    "86-33=53",
    "86+6=92",
    "97-7=90",
    "12+71=83",
    "50-35=15",
    "84-14=70",
    "22+17=39",
    "78-16=62",
    "31+3=34",
    "49-46=3",
    "60-10=50",
    "81-37=44",
    "59+2=61",
    "27-1=26",
    "43+5=48",
    "24+19=43",
    "17+82=99",
    "65+4=69",
    "23-16=7",
    "74-29=45",
    "38-5=33",
    "45-30=15",
    "7+79=86",
    "11+18=29",
    "9+29=38",
    "8+64=72",
    "99-40=59",
    "53+35=88",
    "4+5=9",
    "93-79=14",
    "20-9=11",
    "10+82=92",
    "56-26=30",
    "84-35=49",
    "1+67=68",
    "1+9=10",
    "59+0=59",
    "13+31=44",
    "17+70=87",
    "8+12=20",
    "52+42=94",
    "77-34=43",
    "91-6=85",
    "49+36=85",
    "7+30=37",
    "76-11=65",
    "46-11=35",
    "61-22=39",
    "58+0=58",
    "17+22=39",
    "6+51=57",
    "99-42=57",
    "37+53=90",
    "72-68=4",
    "48+21=69",
    "88-86=2",
    "5+52=57",
    "95-80=15",
    "23+0=23",
    "61-28=33",
    "15-14=1",
    "36+2=38",
    "83+14=97",
    "54+38=92",
    "56-56=0",
    "18+51=69",
    "68-43=25",
    "86-57=29",
    "44+43=87",
    "23+53=76",
    "19-18=1",
    "12+50=62",
    "30+55=85",
    "5+17=22",
    "19+14=33",
    "33-2=31",
    "5+51=56",
    "12+4=16",
    "70+5=75",
    "46+40=86",
    "35+37=72"
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Host "Updated $idx cells"
